$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 3: remove custom height, change D3 font style (16pt bold -> 14pt bold) ---
$ws.Rows.Item(3).RowHeight = 18
$ws.Cells.Item(3,4).Font.Name = "Arial"
$ws.Cells.Item(3,4).Font.Size = 14
$ws.Cells.Item(3,4).Font.Bold = $true
$ws.Cells.Item(3,4).VerticalAlignment = -4108   # xlCenter
$ws.Cells.Item(3,4).WrapText = $true

# --- Row 5 / Row 6: remove custom height (value becomes 18) ---
$ws.Rows.Item(5).RowHeight = 18
$ws.Rows.Item(6).RowHeight = 18

# --- Row 7: D7 loses its "FORMATO X" text, style (s=6) stays ---
$ws.Cells.Item(7,4).ClearContents()

# --- Row 8: now carries the "Fecha:" label (moved up from D9), with a new explicit height ---
$ws.Rows.Item(8).RowHeight = 15.75
$ws.Cells.Item(8,4).Value = "Fecha:"
$ws.Cells.Item(8,4).Font.Name = "Arial"
$ws.Cells.Item(8,4).Font.Size = 12
$ws.Cells.Item(8,4).Font.Bold = $true
$ws.Cells.Item(8,4).HorizontalAlignment = -4152 # xlRight
$ws.Cells.Item(8,4).VerticalAlignment = -4108   # xlCenter
$ws.Cells.Item(8,4).WrapText = $true

# --- Row 9: now carries the "Hora:" label (moved up from D10); keeps its existing style ---
$ws.Cells.Item(9,4).Value = "Hora:"

# --- Row 10: D10 ("Hora:") is removed completely (cell + format) ---
$ws.Cells.Item(10,4).Clear()

# --- Row 12 (becomes row 11 after the blank-row deletion below): "MRO" -> "MRVyO" ---
$ws.Cells.Item(12,5).Value = "MRVyO"

# --- Remove the blank row 11 so row 12's content shifts up to row 11 ---
$ws.Rows.Item(11).Delete()

# --- Update the selected cell shown when the sheet is opened ---
$ws.Range("A6").Select()
